# Sample Project / Main.xlsx : row 11 "From" cell (B11) is changed from the
# text "R40" to the text "1", while keeping its original cell formatting
# (style) untouched.
#
# A plain  $cell.Value = "1"  would make Excel auto-convert the text into a
# number (losing the shared-string/text typing), and forcing text via a
# leading apostrophe (quote-prefix) changes the cell's style (Excel tracks
# the quote-prefix as part of the cell format). To avoid touching the
# style, we stage the quote-prefixed text on a scratch cell, then copy just
# its value back onto B11 so the original style (and everything else) is
# left exactly as it was.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$target = $ws.Cells.Item(11, 2)   # B11
$scratch = $ws.Cells.Item(50, 50) # far-away, unused scratch cell

# Force the scratch cell to hold the literal text "1" (quote-prefix keeps
# Excel from re-interpreting it as a number).
$scratch.Formula = "'1"

# Copy only the computed value (as text) from the scratch cell onto B11,
# leaving B11's existing number format / style untouched.
$scratch.Copy()
$target.PasteSpecial(-4163) | Out-Null  # xlPasteValues

# Clean up the scratch cell so it doesn't leave stray data behind.
$scratch.Clear()
